# Generate Report for Handoff
# Replaces the e2e markdown file id 7eacbddf-... with 2174cc76-...
# and refreshes the associated handoff timestamps / xliff file names.

$wb = $excel.ActiveWorkbook

$oldId = "7eacbddf-b45d-4045-a9d6-3e4957d931f6"
$newId = "2174cc76-3f25-4e35-8aab-b35bd2821a15"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newId.md"
$wsOverview.Range("B2").Value = "e2e\$newId.md"
$wsOverview.Range("G2").Value = "2016-09-02 23:06:46"

# Refresh the hyperlink display text on B2 without touching its target URL.
$overviewUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0d7f953a4986a9c33572b2ce1e80a8f6a170725e/e2e/$oldId.md"
$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $overviewUrl, "", "", "e2e\$newId.md")

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newId.md"
$wsZhCn.Range("G2").Value = "$newId.0dc5567e34b7120ffa73f9373303d5777b8e0f53.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-09-02 23:06:41"
$wsZhCn.Range("I2").Value = ""
$wsZhCn.Range("J2").Value = ""
$wsZhCn.Range("K2").Value = "0001-01-01 00:00:00"

# Drop the now-obsolete "Latest Target File" hyperlink (I2) while keeping A2's.
$zhCnUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0d7f953a4986a9c33572b2ce1e80a8f6a170725e/e2e/$oldId.md"
$wsZhCn.Range("A2").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $zhCnUrl, "", "", "$newId.md")
$wsZhCn.Range("I2").Style = "Normal"
$wsZhCn.Range("J2").Style = "Normal"

$wsZhCn.Columns.Item(9).ColumnWidth = 17.817272004627068
$wsZhCn.Columns.Item(10).ColumnWidth = 20.872143700009268

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newId.md"
$wsDeDe.Range("G2").Value = "$newId.0dc5567e34b7120ffa73f9373303d5777b8e0f53.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-09-02 23:06:46"
$wsDeDe.Range("I2").Value = ""
$wsDeDe.Range("J2").Value = ""
$wsDeDe.Range("K2").Value = "0001-01-01 00:00:00"

# Drop the now-obsolete "Latest Target File" hyperlink (I2) while keeping A2's.
$deDeUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0d7f953a4986a9c33572b2ce1e80a8f6a170725e/e2e/$oldId.md"
$wsDeDe.Range("A2").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $deDeUrl, "", "", "$newId.md")
$wsDeDe.Range("I2").Style = "Normal"
$wsDeDe.Range("J2").Style = "Normal"

$wsDeDe.Columns.Item(9).ColumnWidth = 17.817272004627068
$wsDeDe.Columns.Item(10).ColumnWidth = 20.872143700009268

Write-Host "Applied handback report refresh ($oldId -> $newId)"
